# Update the "Prototyping Study Report" deliverable (row 9 on the
# Checklist sheet) to mark it ready for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# G9 currently holds the "needs review" note; replace it with the
# "Ready for printing" status used elsewhere in the Category column.
$ws.Range("G9").Value = "Ready for printing"

# That note used wrapped text (which is what forced the tall row); the
# new short status doesn't need wrapping, so turn it off and let the
# row height return to the sheet's default.
$ws.Range("G9").WrapText = $false
$ws.Rows.Item(9).AutoFit()

# Leave the selection sitting on C9, where this edit was made.
$ws.Activate()
$ws.Range("C9").Select()

# Restore the theme's dark-text colour to pure black.
$wb.Theme.ThemeColorScheme.Colors(1).RGB = 0
